$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Colors (Excel Font.Color uses BGR-packed long, not RGB hex)
$purple = 10498160   # RGB 7030A0
$blue   = 12611584   # RGB 0070C0
$red    = 255        # RGB FF0000 (already used elsewhere in the sheet -> font id 1)

# --- New notes appended below the table (rows 35-36) -----------------------
# Creating these first establishes the new "purple" font (and reuses the
# existing "red" font) in the same order the source workbook has them.
$c35 = $ws.Range("C35")
$c35.Value = "requires verification--check that it is this way in all the best models"
$c35.Font.Color = $purple
$c35.WrapText = $true

# --- Section 3 wording refinement (drop "usually") --------------------------
$ws.Range("C32").Value = "3.2a. Directions of responses to predictor variables is consistent"

$c36 = $ws.Range("C36")
$c36.Value = "my guess"
$c36.Font.Color = $red
$c36.WrapText = $true

# --- Highlight the "(no/yes)" prediction-supported columns with blue font --
$ws.Range("D6").Font.Color = $blue
$ws.Range("E6").Font.Color = $blue
$ws.Range("G6").Value = "(no/yes)"

$ws.Range("D7").Font.Color = $red
$ws.Range("E7").Font.Color = $red

$ws.Range("D8").Font.Color = $blue
$ws.Range("E8").Font.Color = $blue
$ws.Range("G8").Value = "(no/yes)"

$ws.Range("D19").Font.Color = $blue
$ws.Range("D20").Font.Color = $blue

$d28 = $ws.Range("D28")
$d28.Value = "no"
$d28.Font.Color = $blue

$ws.Range("D30").Font.Color = $blue

$ws.Range("D33").Font.Color = $blue

$ws.Range("D32").Value = "(yes)"

# Zoom + selection to match the saved view state
$ws.Range("D16:D18").Select()
$excel.ActiveWindow.Zoom = 99

# --- Column widths: F and G both become width 8 ----------------------------
$ws.Columns.Item(6).ColumnWidth = 7.16796875
$ws.Columns.Item(7).ColumnWidth = 7.16796875

Write-Output "done"
